$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (shifts existing rows 9-24 down to 10-25)
$ws.Rows.Item(9).Insert()

# New row 9: same record as the (now shifted) row 10, except newer date and lower volume
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value = "Los Lagos"
$ws.Cells.Item(9, 4).Value = 44497
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 300000000
$ws.Cells.Item(9, 7).Value = "Espárragos"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 120
$ws.Cells.Item(9, 11).Value = 1800
$ws.Cells.Item(9, 12).Value = 1800
$ws.Cells.Item(9, 13).Value = 1800
$ws.Cells.Item(9, 14).Value = "$/kilo"
$ws.Cells.Item(9, 15).Value = "Provincia de Linares"
$ws.Cells.Item(9, 16).Value = 1800
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"
